# "Generate Report for Handoff"
#
# The localization-status report was regenerated; the only substantive
# content change is the "Latest Handoff Datetime" for the file
# 45adacbf-760e-419c-8f27-a66d66377ffe.md on the zh-cn handoff sheet,
# which advances from 2017-02-17 07:58:33 to 2017-02-17 07:59:22
# (a new handoff xliff was generated for that language).
#
# (All other cell-index churn visible in the raw OOXML diff is just the
# shared-string table shifting to accommodate the newly-appended date
# string - the same text being reused everywhere else - so no other
# cell's displayed value actually changes.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("zh-cn")

# Row 4 corresponds to 45adacbf-760e-419c-8f27-a66d66377ffe.md; column H is
# "Latest Handoff Datetime". Assign as text so it stays a string cell (style
# already carries the date display format) rather than becoming a date serial.
$ws.Range("H4").Value = "2017-02-17 07:59:22"
